{"js": "// Update the date line and all the division problems in the table to the\n// new values from the latest generated worksheet.\n\nconst replacements = [\n  { find: \"2025-05-26 Monday\", replace: \"2025-05-27 Tuesday\" },\n  { find: \"979\u00f77=\", replace: \"999\u00f73=\" },\n  { find: \"346\u00f79=\", replace: \"337\u00f79=\" },\n  { find: \"293\u00f76=\", replace: \"308\u00f72=\" },\n  { find: \"219\u00f76=\", replace: \"535\u00f76=\" },\n  { find: \"615\u00f74=\", replace: \"198\u00f76=\" },\n  { find: \"746\u00f74=\", replace: \"579\u00f79=\" },\n  { find: \"664\u00f79=\", replace: \"724\u00f79=\" },\n  { find: \"923\u00f77=\", replace: \"429\u00f76=\" },\n  { find: \"615\u00f73=\", replace: \"459\u00f77=\" },\n  { find: \"276\u00f74=\", replace: \"766\u00f79=\" },\n  { find: \"403\u00f74=\", replace: \"322\u00f73=\" },\n  { find: \"695\u00f73=\", replace: \"311\u00f73=\" },\n  { find: \"637\u00f72=\", replace: \"388\u00f79=\" },\n  { find: \"766\u00f75=\", replace: \"309\u00f78=\" },\n  { find: \"153\u00f73=\", replace: \"599\u00f78=\" },\n  { find: \"549\u00f76=\", replace: \"982\u00f74=\" },\n  { find: \"651\u00f72=\", replace: \"877\u00f77=\" },\n  { find: \"104\u00f72=\", replace: \"972\u00f74=\" },\n  { find: \"496\u00f78=\", replace: \"622\u00f78=\" },\n  { find: \"633\u00f77=\", replace: \"843\u00f73=\" },\n  { find: \"594\u00f78=\", replace: \"653\u00f79=\" },\n  { find: \"376\u00f78=\", replace: \"965\u00f74=\" },\n  { find: \"901\u00f78=\", replace: \"983\u00f73=\" },\n  { find: \"262\u00f78=\", replace: \"476\u00f73=\" },\n  { find: \"264\u00f76=\", replace: \"628\u00f78=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all the division problems in the table to the\n# new values from the latest generated worksheet.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-05-26 Monday\"; Replace = \"2025-05-27 Tuesday\" },\n    @{ Find = \"979\u00f77=\"; Replace = \"999\u00f73=\" },\n    @{ Find = \"346\u00f79=\"; Replace = \"337\u00f79=\" },\n    @{ Find = \"293\u00f76=\"; Replace = \"308\u00f72=\" },\n    @{ Find = \"219\u00f76=\"; Replace = \"535\u00f76=\" },\n    @{ Find = \"615\u00f74=\"; Replace = \"198\u00f76=\" },\n    @{ Find = \"746\u00f74=\"; Replace = \"579\u00f79=\" },\n    @{ Find = \"664\u00f79=\"; Replace = \"724\u00f79=\" },\n    @{ Find = \"923\u00f77=\"; Replace = \"429\u00f76=\" },\n    @{ Find = \"615\u00f73=\"; Replace = \"459\u00f77=\" },\n    @{ Find = \"276\u00f74=\"; Replace = \"766\u00f79=\" },\n    @{ Find = \"403\u00f74=\"; Replace = \"322\u00f73=\" },\n    @{ Find = \"695\u00f73=\"; Replace = \"311\u00f73=\" },\n    @{ Find = \"637\u00f72=\"; Replace = \"388\u00f79=\" },\n    @{ Find = \"766\u00f75=\"; Replace = \"309\u00f78=\" },\n    @{ Find = \"153\u00f73=\"; Replace = \"599\u00f78=\" },\n    @{ Find = \"549\u00f76=\"; Replace = \"982\u00f74=\" },\n    @{ Find = \"651\u00f72=\"; Replace = \"877\u00f77=\" },\n    @{ Find = \"104\u00f72=\"; Replace = \"972\u00f74=\" },\n    @{ Find = \"496\u00f78=\"; Replace = \"622\u00f78=\" },\n    @{ Find = \"633\u00f77=\"; Replace = \"843\u00f73=\" },\n    @{ Find = \"594\u00f78=\"; Replace = \"653\u00f79=\" },\n    @{ Find = \"376\u00f78=\"; Replace = \"965\u00f74=\" },\n    @{ Find = \"901\u00f78=\"; Replace = \"983\u00f73=\" },\n    @{ Find = \"262\u00f78=\"; Replace = \"476\u00f73=\" },\n    @{ Find = \"264\u00f76=\"; Replace = \"628\u00f78=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
